$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = -12.339
$ws.Range("B12").Value = 4.935
$ws.Range("C23").Value = -12.609
$ws.Range("D24").Value = -7.606999999999999
$ws.Range("C28").Value = -12.604
$ws.Range("B32").Value = 6.408000000000001
$ws.Range("C32").Value = -12.327
$ws.Range("C34").Value = -11.63
$ws.Range("B36").Value = 8.620000000000001
$ws.Range("B38").Value = 5.547
$ws.Range("D38").Value = -7.802
$ws.Range("C42").Value = -12.221
$ws.Range("B46").Value = 6.248
$ws.Range("D52").Value = -7.780999999999999
$ws.Range("B54").Value = 5.401
$ws.Range("C54").Value = -13.017
$ws.Range("B55").Value = 4.763
$ws.Range("B67").Value = 5.526
$ws.Range("B69").Value = 5.305
$ws.Range("B72").Value = 5.697
$ws.Range("D78").Value = -8.285
$ws.Range("D83").Value = -8.061000000000002
$ws.Range("D85").Value = -8.450000000000001
$ws.Range("D86").Value = -8.354000000000003
$ws.Range("B91").Value = 5.296000000000001
$ws.Range("D96").Value = -7.202000000000001
$ws.Range("C97").Value = -11.479
$ws.Range("B99").Value = 5.524000000000001
$ws.Range("C99").Value = -12.039
$ws.Range("C101").Value = -12.303
$ws.Range("D103").Value = -8.309000000000001
$ws.Range("B104").Value = 8.032
